# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker table (rows 16-19) was regrouped: previously it listed each
# worker together with both of their overdue periods (Iris 1604/1603, then
# Katerine 1604/1603). Now the rows are grouped by period instead, so each
# period (1603, then 1604) lists both workers for that period.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16: Iris Yaneth Rojas Diaz, period 1603 (was 1604)
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "37685377"
$ws.Range("D16").Value = "IRIS YANETH ROJAS DIAZ"
$ws.Range("E16").Value = "1603"

# Row 17: Katerine Yaneth Brieva Rojas, period 1603
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1143388296"
$ws.Range("D17").Value = "KATERINE YANETH BRIEVA ROJAS"
$ws.Range("E17").Value = "1603"

# Row 18: Iris Yaneth Rojas Diaz, period 1604
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "37685377"
$ws.Range("D18").Value = "IRIS YANETH ROJAS DIAZ"
$ws.Range("E18").Value = "1604"

# Row 19: Katerine Yaneth Brieva Rojas, period 1604 (unchanged content)
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1143388296"
$ws.Range("D19").Value = "KATERINE YANETH BRIEVA ROJAS"
$ws.Range("E19").Value = "1604"
